$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.24
$ws.Range("C2").Value = 0.48
$ws.Range("J2").Value = 0.01666666666666667
$ws.Range("P2").Value = 0.1566666666666667
$ws.Range("S2").Value = 0.1066666666666667
$ws.Range("B3").Value = 0.0273972602739726
$ws.Range("C3").Value = 0.0273972602739726
$ws.Range("J3").Value = 0.04794520547945205
$ws.Range("P3").Value = 0.636986301369863
$ws.Range("S3").Value = 0.2602739726027397
$ws.Range("P4").Value = 0.6451612903225806
$ws.Range("S4").Value = 0.3548387096774194
$ws.Range("B6").Value = 0.05263157894736842
$ws.Range("F6").Value = 0.05701754385964912
$ws.Range("J6").Value = 0.2631578947368421
$ws.Range("O6").Value = 0.004385964912280702
$ws.Range("Q6").Value = 0.1271929824561404
$ws.Range("R6").Value = 0.1008771929824561
$ws.Range("S6").Value = 0.3947368421052632
$ws.Range("B7").Value = 0.131578947368421
$ws.Range("D7").Value = 0.03947368421052631
$ws.Range("F7").Value = 0.03947368421052631
$ws.Range("J7").Value = 0.125
$ws.Range("O7").Value = 0.0131578947368421
$ws.Range("Q7").Value = 0.131578947368421
$ws.Range("R7").Value = 0.07894736842105263
$ws.Range("S7").Value = 0.4407894736842105
$ws.Range("B8").Value = 0.09615384615384616
$ws.Range("D8").Value = 0.01201923076923077
$ws.Range("E8").Value = 0.002403846153846154
$ws.Range("F8").Value = 0.0576923076923077
$ws.Range("J8").Value = 0.125
$ws.Range("O8").Value = 0.009615384615384616
$ws.Range("Q8").Value = 0.1706730769230769
$ws.Range("R8").Value = 0.1057692307692308
$ws.Range("S8").Value = 0.4206730769230769
$ws.Range("B9").Value = 0.1120331950207469
$ws.Range("D9").Value = 0.01244813278008299
$ws.Range("F9").Value = 0.06224066390041494
$ws.Range("J9").Value = 0.1037344398340249
$ws.Range("O9").Value = 0.02074688796680498
$ws.Range("Q9").Value = 0.1701244813278008
$ws.Range("R9").Value = 0.09958506224066389
$ws.Range("S9").Value = 0.4190871369294606
$ws.Range("B10").Value = 0.1106194690265487
$ws.Range("D10").Value = 0.01592920353982301
$ws.Range("F10").Value = 0.09557522123893805
$ws.Range("J10").Value = 0.1247787610619469
$ws.Range("O10").Value = 0.0168141592920354
$ws.Range("Q10").Value = 0.2017699115044248
$ws.Range("R10").Value = 0.08849557522123894
$ws.Range("S10").Value = 0.3460176991150443
$ws.Range("G11").Value = 0.1222707423580786
$ws.Range("J11").Value = 0.06986899563318777
$ws.Range("K11").Value = 0.2008733624454148
$ws.Range("L11").Value = 0.5851528384279476
$ws.Range("S11").Value = 0.02183406113537118
$ws.Range("G12").Value = 0.7338129496402878
$ws.Range("J12").Value = 0.223021582733813
$ws.Range("L12").Value = 0.02877697841726619
$ws.Range("S12").Value = 0.01438848920863309
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.2857142857142857
$ws.Range("S13").Value = 0.04761904761904762
$ws.Range("F15").Value = 0.02604166666666667
$ws.Range("H15").Value = 0.140625
$ws.Range("I15").Value = 0.078125
$ws.Range("J15").Value = 0.375
$ws.Range("K15").Value = 0.07291666666666667
$ws.Range("M15").Value = 0.015625
$ws.Range("O15").Value = 0.09895833333333333
$ws.Range("S15").Value = 0.1927083333333333
$ws.Range("F16").Value = 0.02597402597402598
$ws.Range("H16").Value = 0.1493506493506493
$ws.Range("I16").Value = 0.1233766233766234
$ws.Range("J16").Value = 0.3571428571428572
$ws.Range("K16").Value = 0.09740259740259741
$ws.Range("M16").Value = 0.06493506493506493
$ws.Range("O16").Value = 0.04545454545454546
$ws.Range("S16").Value = 0.1363636363636364
$ws.Range("F17").Value = 0.01813471502590673
$ws.Range("H17").Value = 0.1917098445595855
$ws.Range("I17").Value = 0.1295336787564767
$ws.Range("J17").Value = 0.4119170984455959
$ws.Range("K17").Value = 0.07253886010362694
$ws.Range("M17").Value = 0.01036269430051814
$ws.Range("O17").Value = 0.06994818652849741
$ws.Range("S17").Value = 0.09585492227979274
$ws.Range("F18").Value = 0.02450980392156863
$ws.Range("H18").Value = 0.142156862745098
$ws.Range("I18").Value = 0.142156862745098
$ws.Range("J18").Value = 0.4264705882352941
$ws.Range("K18").Value = 0.1176470588235294
$ws.Range("M18").Value = 0.0196078431372549
$ws.Range("S18").Value = 0.06862745098039216
$ws.Range("F19").Value = 0.01650738488271069
$ws.Range("H19").Value = 0.2258905299739357
$ws.Range("I19").Value = 0.1112076455256299
$ws.Range("J19").Value = 0.3544743701129452
$ws.Range("K19").Value = 0.08601216333622937
$ws.Range("M19").Value = 0.01824500434404865
$ws.Range("N19").Value = 0.0008688097306689834
$ws.Range("O19").Value = 0.06602953953084274
$ws.Range("S19").Value = 0.1207645525629887
